$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The guest data row stores every value (including numeric-looking and
# date-looking ones) as plain text/shared strings. Force text formatting
# before assigning so Excel does not auto-convert "1"/"2" to numbers or
# "2025-01-11" to a date serial, then clear the temporary formatting so
# no residual style/number-format is left on the cells.
$cellA2 = $ws.Range("A2")
$cellB2 = $ws.Range("B2")
$cellD2 = $ws.Range("D2")
$cellE2 = $ws.Range("E2")

$cellA2.NumberFormat = "@"
$cellB2.NumberFormat = "@"
$cellD2.NumberFormat = "@"
$cellE2.NumberFormat = "@"

$cellA2.Value = "Shreed"
$cellB2.Value = "1"
$cellD2.Value = "2025-01-11"
$cellE2.Value = "2"

$cellA2.ClearFormats()
$cellB2.ClearFormats()
$cellD2.ClearFormats()
$cellE2.ClearFormats()
